$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")
$wsAssets    = $wb.Worksheets.Item("Assets")

# ---------------------------------------------------------------------------
# Settings sheet
# ---------------------------------------------------------------------------
# Row 2: OrchestratorQueueName / ProcessABCQueue / Orchestrator queue Name...
#   -> InputDataFilePath / Data\Input\Question3.xlsx / Place where the input
#      file is stored
$wsSettings.Range("A2").Value = "InputDataFilePath"
$wsSettings.Range("B2").Value = "Data\Input\Question3.xlsx"
$wsSettings.Range("C2").Value = "Place where the input file is stored"

# Row 3: OrchestratorQueueFolder / (empty) / Folder name...
#   -> InputDataSheetName / Notaries Area of Practice / Sheet name where the
#      input data is stored
$wsSettings.Range("A3").Value = "InputDataSheetName"
$wsSettings.Range("B3").Value = "Notaries Area of Practice"
$wsSettings.Range("C3").Value = "Sheet name where the input data is stored"
# The old row-3 text wrapped across 3 lines; the new text fits on one, so the
# custom wrap height no longer applies - let the row shrink back to normal.
$wsSettings.Rows.Item(3).AutoFit()

# Row 5: logF_BusinessProcessName / Framework / Logging field...
#   -> logF_BusinessProcessName / BotsDna - Notaries / Logging field... (text unchanged)
$wsSettings.Range("B5").Value = "BotsDna – Notaries"
$wsSettings.Rows.Item(5).RowHeight = 30

# Row 6: new content OutputSheet / Sheet2
$wsSettings.Range("A6").Value = "OutputSheet"
$wsSettings.Range("B6").Value = "Sheet2"

# ---------------------------------------------------------------------------
# Constants sheet
# ---------------------------------------------------------------------------
# Row 18: new content BrowserPath / https://botsdna.com/notaries/
$wsConstants.Range("A18").Value = "BrowserPath"
$wsConstants.Range("B18").Value = "https://botsdna.com/notaries/"

# Update the stored selection on Constants (A12 -> A19) without leaving it
# as the active sheet - the Settings sheet becomes active/selected last.
$wsConstants.Range("A19").Select() | Out-Null

# ---------------------------------------------------------------------------
# Assets sheet - no content changes; it merely stops being the active tab.
# ---------------------------------------------------------------------------

# Settings becomes the active sheet/tab, with B6 as the selected cell -
# activate and select last so it "wins" as the workbook's active sheet.
$wsSettings.Activate()
$wsSettings.Range("B6").Select() | Out-Null
